$d = $word.ActiveDocument
$brk = [char]11

# --- Change 1: remove the "LOB1009 - Leitura e Interpretação de Desenho
#     Técnico (Requisito)" line entirely (whole run incl. its line break) ---
$find1 = $d.Content.Find
$find1.ClearFormatting()
$find1.Text = "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)" + $brk
$find1.Replacement.ClearFormatting()
$find1.Replacement.Text = ""
$ok1 = $find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2)

# --- Change 2: remove the three existing Química requirement lines
#     (LOQ4031, LOQ4073, LOQ4095) as a single block so the rest of the
#     paragraph / runs are left untouched ---
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Text = "LOQ4031 -  Química Geral I  (Requisito)" + $brk + "LOQ4073 -  Química Geral II  (Requisito)" + $brk + "LOQ4095 -  Química Geral Experimental  (Requisito)" + $brk
$find2.Replacement.ClearFormatting()
$find2.Replacement.Text = ""
$ok2 = $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)

# --- Change 3: append the four replacement Química lines, one run at a
#     time, at the very end of the (last) paragraph. Inserting after a
#     collapsed range at the true end of the paragraph (rather than in
#     the middle of existing runs, and rather than InsertBefore) is what
#     makes each addition land as its own separate <w:r> run instead of
#     being folded into a neighbouring run, matching the target shape ---
$line1 = "LOQ4095 -  Química Geral Experimental  (Requisito)" + $brk
$p1 = $d.Paragraphs.Last
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertAfter($line1)

$line2 = "LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito)" + $brk
$p2 = $d.Paragraphs.Last
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertAfter($line2)

$line3 = "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)" + $brk
$p3 = $d.Paragraphs.Last
$r3 = $p3.Range
$r3.Collapse(0)
$r3.InsertAfter($line3)

$line4 = "LOQ4247 -  Desenho Assistido por Computador  (Requisito)" + $brk
$p4 = $d.Paragraphs.Last
$r4 = $p4.Range
$r4.Collapse(0)
$r4.InsertAfter($line4)

Write-Output "LOB1009 removed: $ok1; LOQ block removed: $ok2; 4 new LOQ lines appended"
